$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) updates are written as literal text, matching the
# source data (t="inlineStr") rather than being auto-parsed as numbers/dates.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.288.64'
$ws.Range("E2").Value = '  +2.48%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.05'
$ws.Range("E3").Value = '  +0.00%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.17'
$ws.Range("E5").Value = '  +3.59%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.51'
$ws.Range("E6").Value = '  +1.90%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("E7").Value = '  +0.53%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +5.36%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.44'
$ws.Range("E10").Value = '  +0.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0799'
$ws.Range("E11").Value = '  -0.14%  '

# Row 12
$ws.Range("E12").Value = '  -1.59%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.30'
$ws.Range("E13").Value = '  -2.77%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.04'
$ws.Range("E14").Value = '  +1.67%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.804.62'
$ws.Range("E15").Value = '  +0.07%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.415.33'
$ws.Range("E16").Value = '  +0.05%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("E17").Value = '  +1.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.187.11'
$ws.Range("E18").Value = '  +2.38%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.21'
$ws.Range("E19").Value = '  -0.66%  '

# Row 20
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("E21").Value = '  +1.71%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.99'
$ws.Range("E22").Value = '  +0.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.72'
$ws.Range("E23").Value = '  +1.89%  '

# Row 24
$ws.Range("E24").Value = '  -0.69%  '

# Row 25
$ws.Range("E25").Value = '  +0.15%  '

# Row 26
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.74'
$ws.Range("E27").Value = '  +2.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  -2.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.62'
$ws.Range("E29").Value = '  +0.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.57'
$ws.Range("E30").Value = '  +2.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.95'
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.16'
$ws.Range("E32").Value = '  +8.23%  '

# Row 33
$ws.Range("E33").Value = '  +7.93%  '

# Row 34
$ws.Range("E34").Value = '  +0.88%  '

# Row 35
$ws.Range("E35").Value = '  +0.19%  '

# Row 36
$ws.Range("E36").Value = '  +0.39%  '

# Row 37
$ws.Range("E37").Value = '  -1.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.44'
$ws.Range("E38").Value = '  -0.77%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '127.00'
$ws.Range("E39").Value = '  -1.96%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("E40").Value = '  -0.31%  '

# Row 41
$ws.Range("E41").Value = '  +0.60%  '

# Row 42
$ws.Range("E42").Value = '  -3.95%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.72'
$ws.Range("E43").Value = '  -5.28%  '

# Row 44
$ws.Range("E44").Value = '  +0.77%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.937.67'
$ws.Range("E45").Value = '  -0.92%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.12'
$ws.Range("E48").Value = '  -2.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.79'
$ws.Range("E49").Value = '  +8.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.46'
$ws.Range("E50").Value = '  +3.72%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.84'
$ws.Range("E51").Value = '  +6.31%  '
